# Auto-generated script to apply scheduled market-data refresh to the Leve profit workbook.
# Updates columns H-N (current market prices / computed profit figures) across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the latest Universalis price snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 948.75
$ws.Range("I28").Value = 831.53845
$ws.Range("K28").Value = 831.53845
$ws.Range("M28").Value = -346.53845
$ws.Range("H31").Value = 5926.6665
$ws.Range("I31").Value = 5945
$ws.Range("J31").Value = 5890
$ws.Range("K31").Value = 17835
$ws.Range("L31").Value = 17670
$ws.Range("M31").Value = -17605
$ws.Range("N31").Value = -18130
$ws.Range("H33").Value = 283.0435
$ws.Range("J33").Value = 525
$ws.Range("L33").Value = 525
$ws.Range("N33").Value = -983
$ws.Range("H62").Value = 7269.75
$ws.Range("I62").Value = 7814.6665
$ws.Range("J62").Value = 6569.143
$ws.Range("K62").Value = 7814.6665
$ws.Range("L62").Value = 6569.143
$ws.Range("M62").Value = -7190.6665
$ws.Range("N62").Value = -7817.143
$ws.Range("H64").Value = 22037.084
$ws.Range("I64").Value = 26938.334
$ws.Range("J64").Value = 7333.3335
$ws.Range("K64").Value = 26938.334
$ws.Range("L64").Value = 7333.3335
$ws.Range("M64").Value = -26690.334
$ws.Range("N64").Value = -7829.3335
$ws.Range("H65").Value = 7269.75
$ws.Range("I65").Value = 7814.6665
$ws.Range("J65").Value = 6569.143
$ws.Range("K65").Value = 39073.3325
$ws.Range("L65").Value = 32845.715
$ws.Range("M65").Value = -35953.3325
$ws.Range("N65").Value = -39085.715
$ws.Range("H67").Value = 22037.084
$ws.Range("I67").Value = 26938.334
$ws.Range("J67").Value = 7333.3335
$ws.Range("K67").Value = 26938.334
$ws.Range("L67").Value = 7333.3335
$ws.Range("M67").Value = -26080.334
$ws.Range("N67").Value = -9049.333500000001
$ws.Range("H80").Value = 2180338
$ws.Range("I80").Value = 4831629
$ws.Range("J80").Value = 11100.272
$ws.Range("K80").Value = 14494887
$ws.Range("L80").Value = 33300.81600000001
$ws.Range("M80").Value = -14493889
$ws.Range("N80").Value = -35296.81600000001
$ws.Range("H83").Value = 2180338
$ws.Range("I83").Value = 4831629
$ws.Range("J83").Value = 11100.272
$ws.Range("K83").Value = 43484661
$ws.Range("L83").Value = 99902.448
$ws.Range("M83").Value = -43479669
$ws.Range("N83").Value = -109886.448
$ws.Range("H100").Value = 2887.375
$ws.Range("I100").Value = 1350.3334
$ws.Range("J100").Value = 7498.5
$ws.Range("K100").Value = 1350.3334
$ws.Range("L100").Value = 7498.5
$ws.Range("M100").Value = -809.3334
$ws.Range("N100").Value = -8580.5
$ws.Range("H112").Value = 5518.147
$ws.Range("J112").Value = 5852.5483
$ws.Range("L112").Value = 17557.6449
$ws.Range("N112").Value = -19773.6449
$ws.Range("H114").Value = 59999
$ws.Range("J114").Value = 59999
$ws.Range("L114").Value = 59999
$ws.Range("N114").Value = -68677
$ws.Range("H116").Value = 847291.5600000001
$ws.Range("I116").Value = 1267543.4
$ws.Range("K116").Value = 1267543.4
$ws.Range("M116").Value = -1264101.4
$ws.Range("H127").Value = 1365.3334
$ws.Range("I127").Value = 1238.4
$ws.Range("J127").Value = 2000
$ws.Range("K127").Value = 3715.2
$ws.Range("L127").Value = 6000
$ws.Range("M127").Value = 1244.8
$ws.Range("N127").Value = -15920
$ws.Range("H132").Value = 3191417.8
$ws.Range("J132").Value = 2617.6667
$ws.Range("L132").Value = 7853.000100000001
$ws.Range("N132").Value = -12913.0001
$ws.Range("H137").Value = 20379.666
$ws.Range("I137").Value = 21372
$ws.Range("J137").Value = 19445.705
$ws.Range("K137").Value = 64116
$ws.Range("L137").Value = 58337.11500000001
$ws.Range("M137").Value = -61566
$ws.Range("N137").Value = -63437.11500000001
$ws.Range("H138").Value = 52274.383
$ws.Range("I138").Value = 3657.1428
$ws.Range("J138").Value = 149508.86
$ws.Range("K138").Value = 10971.4284
$ws.Range("L138").Value = 448526.58
$ws.Range("M138").Value = -5831.428400000001
$ws.Range("N138").Value = -458806.58
$ws.Range("H141").Value = 1560.0667
$ws.Range("I141").Value = 1560.0667
$ws.Range("K141").Value = 4680.2001
$ws.Range("M141").Value = 499.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1287.0834
$ws.Range("I2").Value = 1212.6086
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 1212.6086
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -1099.6086
$ws.Range("N2").Value = -3226
$ws.Range("H32").Value = 16457.041
$ws.Range("I32").Value = 16457.041
$ws.Range("K32").Value = 16457.041
$ws.Range("M32").Value = -16170.041
$ws.Range("H61").Value = 12515.9
$ws.Range("I61").Value = 6337.8125
$ws.Range("J61").Value = 37228.25
$ws.Range("K61").Value = 6337.8125
$ws.Range("L61").Value = 37228.25
$ws.Range("M61").Value = -6125.8125
$ws.Range("N61").Value = -37652.25
$ws.Range("H74").Value = 162295.73
$ws.Range("I74").Value = 223437.67
$ws.Range("J74").Value = 12220.091
$ws.Range("K74").Value = 223437.67
$ws.Range("L74").Value = 12220.091
$ws.Range("M74").Value = -222563.67
$ws.Range("N74").Value = -13968.091
$ws.Range("H77").Value = 162295.73
$ws.Range("I77").Value = 223437.67
$ws.Range("J77").Value = 12220.091
$ws.Range("K77").Value = 1117188.35
$ws.Range("L77").Value = 61100.455
$ws.Range("M77").Value = -1112820.35
$ws.Range("N77").Value = -69836.455
$ws.Range("H110").Value = 39620.824
$ws.Range("I110").Value = 42489.477
$ws.Range("K110").Value = 42489.477
$ws.Range("M110").Value = -40444.477
$ws.Range("H116").Value = 1287.0834
$ws.Range("I116").Value = 1212.6086
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1212.6086
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1081.3914
$ws.Range("N116").Value = -7588
$ws.Range("H122").Value = 1833.931
$ws.Range("I122").Value = 1756.5714
$ws.Range("K122").Value = 5269.7142
$ws.Range("M122").Value = -2819.7142
$ws.Range("H132").Value = 2228.7273
$ws.Range("I132").Value = 2001.898
$ws.Range("J132").Value = 4081.1667
$ws.Range("K132").Value = 6005.694
$ws.Range("L132").Value = 12243.5001
$ws.Range("M132").Value = -3475.694
$ws.Range("N132").Value = -17303.5001
$ws.Range("H136").Value = 12515.9
$ws.Range("I136").Value = 6337.8125
$ws.Range("J136").Value = 37228.25
$ws.Range("K136").Value = 19013.4375
$ws.Range("L136").Value = 111684.75
$ws.Range("M136").Value = -16463.4375
$ws.Range("N136").Value = -116784.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1287.0834
$ws.Range("I3").Value = 1212.6086
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 1212.6086
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -1098.6086
$ws.Range("N3").Value = -3228
$ws.Range("H20").Value = 1159.4333
$ws.Range("I20").Value = 1059.35
$ws.Range("J20").Value = 1359.6
$ws.Range("K20").Value = 1059.35
$ws.Range("L20").Value = 1359.6
$ws.Range("M20").Value = -812.3499999999999
$ws.Range("N20").Value = -1853.6
$ws.Range("H80").Value = 724.5625
$ws.Range("I80").Value = 1053.7142
$ws.Range("K80").Value = 1053.7142
$ws.Range("M80").Value = -55.71419999999989
$ws.Range("H83").Value = 724.5625
$ws.Range("I83").Value = 1053.7142
$ws.Range("K83").Value = 5268.571
$ws.Range("M83").Value = -276.5709999999999
$ws.Range("H86").Value = 2333
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2999
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2999
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -5245
$ws.Range("H89").Value = 2333
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2999
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 14995
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -26227
$ws.Range("H99").Value = 1301.7778
$ws.Range("I99").Value = 1115.862
$ws.Range("J99").Value = 2072
$ws.Range("K99").Value = 1115.862
$ws.Range("L99").Value = 2072
$ws.Range("M99").Value = 382.1379999999999
$ws.Range("N99").Value = -5068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 714.5
$ws.Range("I22").Value = 595.2727
$ws.Range("J22").Value = 815.38464
$ws.Range("K22").Value = 595.2727
$ws.Range("L22").Value = 815.38464
$ws.Range("M22").Value = -245.2727
$ws.Range("N22").Value = -1515.38464
$ws.Range("H31").Value = 1221330.5
$ws.Range("I31").Value = 1853493.9
$ws.Range("J31").Value = 2158.1072
$ws.Range("K31").Value = 1853493.9
$ws.Range("L31").Value = 2158.1072
$ws.Range("M31").Value = -1853198.9
$ws.Range("N31").Value = -2748.1072
$ws.Range("H34").Value = 1221330.5
$ws.Range("I34").Value = 1853493.9
$ws.Range("J34").Value = 2158.1072
$ws.Range("K34").Value = 1853493.9
$ws.Range("L34").Value = 2158.1072
$ws.Range("M34").Value = -1853291.9
$ws.Range("N34").Value = -2562.1072
$ws.Range("H50").Value = 36897.6
$ws.Range("J50").Value = 36897.6
$ws.Range("L50").Value = 36897.6
$ws.Range("N50").Value = -38147.6
$ws.Range("H51").Value = 37000
$ws.Range("J51").Value = 37000
$ws.Range("L51").Value = 37000
$ws.Range("N51").Value = -38472
$ws.Range("H58").Value = 1167.6487
$ws.Range("I58").Value = 942.871
$ws.Range("J58").Value = 2329
$ws.Range("K58").Value = 942.871
$ws.Range("L58").Value = 2329
$ws.Range("M58").Value = -739.871
$ws.Range("N58").Value = -2735
$ws.Range("H60").Value = 30948
$ws.Range("I60").Value = 30900
$ws.Range("J60").Value = 30964
$ws.Range("K60").Value = 30900
$ws.Range("L60").Value = 30964
$ws.Range("M60").Value = -30389
$ws.Range("N60").Value = -31986
$ws.Range("H61").Value = 37000
$ws.Range("J61").Value = 37000
$ws.Range("L61").Value = 37000
$ws.Range("N61").Value = -37696
$ws.Range("H62").Value = 2877.8
$ws.Range("I62").Value = 2877.8
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2877.8
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2253.8
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2877.8
$ws.Range("I65").Value = 2877.8
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14389
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11269
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 83333
$ws.Range("J74").Value = 83333
$ws.Range("L74").Value = 83333
$ws.Range("N74").Value = -85081
$ws.Range("H77").Value = 83333
$ws.Range("J77").Value = 83333
$ws.Range("L77").Value = 249999
$ws.Range("N77").Value = -258735
$ws.Range("H86").Value = 13660.7
$ws.Range("I86").Value = 13810.583
$ws.Range("K86").Value = 13810.583
$ws.Range("M86").Value = -12687.583
$ws.Range("H89").Value = 13660.7
$ws.Range("I89").Value = 13810.583
$ws.Range("K89").Value = 69052.91500000001
$ws.Range("M89").Value = -63436.91500000001
$ws.Range("H99").Value = 5860.3076
$ws.Range("I99").Value = 4954.4443
$ws.Range("K99").Value = 4954.4443
$ws.Range("M99").Value = -3456.4443
$ws.Range("H105").Value = 1984.1538
$ws.Range("I105").Value = 1387.875
$ws.Range("K105").Value = 1387.875
$ws.Range("M105").Value = 359.125
$ws.Range("H107").Value = 766.4706
$ws.Range("I107").Value = 571.6923
$ws.Range("J107").Value = 1399.5
$ws.Range("K107").Value = 571.6923
$ws.Range("L107").Value = 1399.5
$ws.Range("M107").Value = 1348.3077
$ws.Range("N107").Value = -5239.5
$ws.Range("H126").Value = 5860.3076
$ws.Range("I126").Value = 4954.4443
$ws.Range("K126").Value = 14863.3329
$ws.Range("M126").Value = -12393.3329
$ws.Range("H132").Value = 44596.652
$ws.Range("I132").Value = 63242.625
$ws.Range("K132").Value = 189727.875
$ws.Range("M132").Value = -187197.875
$ws.Range("H136").Value = 1167.6487
$ws.Range("I136").Value = 942.871
$ws.Range("J136").Value = 2329
$ws.Range("K136").Value = 2828.613
$ws.Range("L136").Value = 6987
$ws.Range("M136").Value = -278.6129999999998
$ws.Range("N136").Value = -12087

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 550.4761999999999
$ws.Range("J5").Value = 574.7
$ws.Range("L5").Value = 1724.1
$ws.Range("N5").Value = -1948.1
$ws.Range("H68").Value = 3747.0364
$ws.Range("I68").Value = 1914.25
$ws.Range("J68").Value = 4258.5117
$ws.Range("K68").Value = 5742.75
$ws.Range("L68").Value = 12775.5351
$ws.Range("M68").Value = -4931.75
$ws.Range("N68").Value = -14397.5351
$ws.Range("H69").Value = 4463.56
$ws.Range("I69").Value = 1496.6666
$ws.Range("J69").Value = 4868.136
$ws.Range("K69").Value = 4489.9998
$ws.Range("L69").Value = 14604.408
$ws.Range("M69").Value = -3678.9998
$ws.Range("N69").Value = -16226.408
$ws.Range("H71").Value = 3747.0364
$ws.Range("I71").Value = 1914.25
$ws.Range("J71").Value = 4258.5117
$ws.Range("K71").Value = 17228.25
$ws.Range("L71").Value = 38326.6053
$ws.Range("M71").Value = -13172.25
$ws.Range("N71").Value = -46438.6053
$ws.Range("H72").Value = 4463.56
$ws.Range("I72").Value = 1496.6666
$ws.Range("J72").Value = 4868.136
$ws.Range("K72").Value = 13469.9994
$ws.Range("L72").Value = 43813.224
$ws.Range("M72").Value = -9413.999400000001
$ws.Range("N72").Value = -51925.224
$ws.Range("H74").Value = 5000
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -17122
$ws.Range("H77").Value = 5000
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -55608
$ws.Range("H107").Value = 2811.125
$ws.Range("I107").Value = 3142.7144
$ws.Range("J107").Value = 490
$ws.Range("K107").Value = 9428.143199999999
$ws.Range("L107").Value = 1470
$ws.Range("M107").Value = -7508.143199999999
$ws.Range("N107").Value = -5310
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 975
$ws.Range("K122").Value = 10800
$ws.Range("L122").Value = 8775
$ws.Range("M122").Value = -8350
$ws.Range("N122").Value = -13675
$ws.Range("H131").Value = 2289.6416
$ws.Range("I131").Value = 3728.1667
$ws.Range("J131").Value = 1868.6097
$ws.Range("K131").Value = 11184.5001
$ws.Range("L131").Value = 5605.8291
$ws.Range("M131").Value = -6144.500100000001
$ws.Range("N131").Value = -15685.8291
$ws.Range("H135").Value = 550.4761999999999
$ws.Range("J135").Value = 574.7
$ws.Range("L135").Value = 5172.3
$ws.Range("N135").Value = -10242.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 474.66666
$ws.Range("I2").Value = 603.6
$ws.Range("J2").Value = 259.77777
$ws.Range("K2").Value = 603.6
$ws.Range("L2").Value = 259.77777
$ws.Range("M2").Value = -490.6
$ws.Range("N2").Value = -485.77777
$ws.Range("H80").Value = 18549.75
$ws.Range("I80").Value = 1999
$ws.Range("J80").Value = 24066.666
$ws.Range("K80").Value = 1999
$ws.Range("L80").Value = 24066.666
$ws.Range("M80").Value = -1001
$ws.Range("N80").Value = -26062.666
$ws.Range("H83").Value = 18549.75
$ws.Range("I83").Value = 1999
$ws.Range("J83").Value = 24066.666
$ws.Range("K83").Value = 9995
$ws.Range("L83").Value = 120333.33
$ws.Range("M83").Value = -5003
$ws.Range("N83").Value = -130317.33
$ws.Range("H97").Value = 1147.7354
$ws.Range("I97").Value = 1068.0476
$ws.Range("J97").Value = 1276.4615
$ws.Range("K97").Value = 1068.0476
$ws.Range("L97").Value = 1276.4615
$ws.Range("M97").Value = -572.0476000000001
$ws.Range("N97").Value = -2268.4615
$ws.Range("H113").Value = 1288.2
$ws.Range("J113").Value = 1297.5
$ws.Range("L113").Value = 1297.5
$ws.Range("N113").Value = -5637.5
$ws.Range("H116").Value = 59999
$ws.Range("J116").Value = 59999
$ws.Range("L116").Value = 59999
$ws.Range("N116").Value = -69177
$ws.Range("H122").Value = 13891537
$ws.Range("I122").Value = 2535.1785
$ws.Range("J122").Value = 62503044
$ws.Range("K122").Value = 7605.5355
$ws.Range("L122").Value = 187509132
$ws.Range("M122").Value = -5155.5355
$ws.Range("N122").Value = -187514032
$ws.Range("H126").Value = 2903.923
$ws.Range("I126").Value = 2345.5454
$ws.Range("J126").Value = 5975
$ws.Range("K126").Value = 7036.6362
$ws.Range("L126").Value = 17925
$ws.Range("M126").Value = -4566.6362
$ws.Range("N126").Value = -22865
$ws.Range("H132").Value = 2287.0571
$ws.Range("I132").Value = 2047.7241
$ws.Range("J132").Value = 3443.8333
$ws.Range("K132").Value = 6143.1723
$ws.Range("L132").Value = 10331.4999
$ws.Range("M132").Value = -3613.1723
$ws.Range("N132").Value = -15391.4999
$ws.Range("H139").Value = 109930
$ws.Range("J139").Value = 109930
$ws.Range("L139").Value = 109930
$ws.Range("N139").Value = -120210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3563
$ws.Range("I7").Value = 4100.8
$ws.Range("K7").Value = 4100.8
$ws.Range("M7").Value = -3988.8
$ws.Range("H22").Value = 1429
$ws.Range("I22").Value = 1044.8889
$ws.Range("J22").Value = 1743.2727
$ws.Range("K22").Value = 1044.8889
$ws.Range("L22").Value = 1743.2727
$ws.Range("M22").Value = -749.8888999999999
$ws.Range("N22").Value = -2333.2727
$ws.Range("H27").Value = 1429
$ws.Range("I27").Value = 1044.8889
$ws.Range("J27").Value = 1743.2727
$ws.Range("K27").Value = 1044.8889
$ws.Range("L27").Value = 1743.2727
$ws.Range("M27").Value = -937.8888999999999
$ws.Range("N27").Value = -1957.2727
$ws.Range("H40").Value = 4082.4167
$ws.Range("I40").Value = 3898.9
$ws.Range("K40").Value = 3898.9
$ws.Range("M40").Value = -3762.9
$ws.Range("H61").Value = 1485.9
$ws.Range("I61").Value = 1252.8572
$ws.Range("K61").Value = 1252.8572
$ws.Range("M61").Value = -1050.8572
$ws.Range("H68").Value = 2139.1853
$ws.Range("I68").Value = 1386.4615
$ws.Range("J68").Value = 2838.1428
$ws.Range("K68").Value = 1386.4615
$ws.Range("L68").Value = 2838.1428
$ws.Range("M68").Value = -637.4614999999999
$ws.Range("N68").Value = -4336.1428
$ws.Range("H71").Value = 2139.1853
$ws.Range("I71").Value = 1386.4615
$ws.Range("J71").Value = 2838.1428
$ws.Range("K71").Value = 6932.307499999999
$ws.Range("L71").Value = 14190.714
$ws.Range("M71").Value = -3188.307499999999
$ws.Range("N71").Value = -21678.714
$ws.Range("H82").Value = 1047.5333
$ws.Range("I82").Value = 1278.6
$ws.Range("J82").Value = 932
$ws.Range("K82").Value = 1278.6
$ws.Range("L82").Value = 932
$ws.Range("M82").Value = -917.5999999999999
$ws.Range("N82").Value = -1654
$ws.Range("H85").Value = 1047.5333
$ws.Range("I85").Value = 1278.6
$ws.Range("J85").Value = 932
$ws.Range("K85").Value = 1278.6
$ws.Range("L85").Value = 932
$ws.Range("M85").Value = -30.59999999999991
$ws.Range("N85").Value = -3428
$ws.Range("H109").Value = 69092.5
$ws.Range("J109").Value = 69092.5
$ws.Range("L109").Value = 69092.5
$ws.Range("N109").Value = -71866.5
$ws.Range("H113").Value = 1485.9
$ws.Range("I113").Value = 1252.8572
$ws.Range("K113").Value = 1252.8572
$ws.Range("M113").Value = 917.1428000000001
$ws.Range("H123").Value = 90000
$ws.Range("J123").Value = 90000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -99800
$ws.Range("H126").Value = 3563
$ws.Range("I126").Value = 4100.8
$ws.Range("K126").Value = 12302.4
$ws.Range("M126").Value = -9832.400000000001
$ws.Range("H132").Value = 2935.2307
$ws.Range("I132").Value = 2461.15
$ws.Range("K132").Value = 7383.450000000001
$ws.Range("M132").Value = -4853.450000000001
$ws.Range("H136").Value = 4011.36
$ws.Range("I136").Value = 4049.7
$ws.Range("K136").Value = 12149.1
$ws.Range("M136").Value = -9599.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 84694.5
$ws.Range("J93").Value = 84694.5
$ws.Range("L93").Value = 84694.5
$ws.Range("N93").Value = -89686.5
$ws.Range("H107").Value = 1100.9524
$ws.Range("I107").Value = 1042.3125
$ws.Range("J107").Value = 1288.6
$ws.Range("K107").Value = 3126.9375
$ws.Range("L107").Value = 3865.8
$ws.Range("M107").Value = -1206.9375
$ws.Range("N107").Value = -7705.799999999999
$ws.Range("H126").Value = 1253550
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 2502100
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 7506300
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -7511240
$ws.Range("H132").Value = 31456.64
$ws.Range("I132").Value = 40326.895
$ws.Range("J132").Value = 3367.5
$ws.Range("K132").Value = 120980.685
$ws.Range("L132").Value = 10102.5
$ws.Range("M132").Value = -118450.685
$ws.Range("N132").Value = -15162.5
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
$ws.Range("H136").Value = 25251.893
$ws.Range("I136").Value = 26001.963
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 78005.889
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -75455.889
$ws.Range("N136").Value = -20100
